# DB Sample case checkin
# Add two new test-scenario rows to the "Scenarios" sheet describing the
# new DB validation sample test class.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# Row 13: dbCancelReasonVerify (not executed)
$ws.Range("A13").Value = "N"
$ws.Range("B13").Value = "patient.web.tests.DBValidationSample"
$ws.Range("D13").Value = "dbCancelReasonVerify"

# Row 14: dbUserAccountVerify (not executed)
$ws.Range("A14").Value = "N"
$ws.Range("B14").Value = "patient.web.tests.DBValidationSample"
$ws.Range("D14").Value = "dbUserAccountVerify"

# Update the sheet view/selection to match the post-edit state.
[void]$ws.Range("A14").Select()

$pws = $wb.Worksheets.Item("Parameters")
[void]$pws.Range("B2").Select()

[void]$ws.Activate()
